# "Actualización 11 de Mayo - Tarde"
# On the "Rescatables" sheet, the two students in rows 3 and 4
# (DE LOS SANTOS GONZALEZ MARIA FERNANDA and GIL BANDALA AELEN) swap
# places, and the "Reprobadas" value that ends up on row 4
# (DE LOS SANTOS ...) becomes 1 instead of 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Capture row 3 (NC, Paterno, Materno, Nombres, Nombre_Largo, Grupo, Reprobadas)
# Use .Value2 so we get plain scalars back, not live COM property wrappers.
$row3 = @(
    $ws.Cells.Item(3, 1).Value2,
    $ws.Cells.Item(3, 2).Value2,
    $ws.Cells.Item(3, 3).Value2,
    $ws.Cells.Item(3, 4).Value2,
    $ws.Cells.Item(3, 5).Value2,
    $ws.Cells.Item(3, 6).Value2,
    $ws.Cells.Item(3, 7).Value2
)

# Capture row 4
$row4 = @(
    $ws.Cells.Item(4, 1).Value2,
    $ws.Cells.Item(4, 2).Value2,
    $ws.Cells.Item(4, 3).Value2,
    $ws.Cells.Item(4, 4).Value2,
    $ws.Cells.Item(4, 5).Value2,
    $ws.Cells.Item(4, 6).Value2,
    $ws.Cells.Item(4, 7).Value2
)

# Write old row 4 data into row 3, keeping "Reprobadas" as it was (2)
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(3, $c).Value = $row4[$c - 1]
}
$ws.Cells.Item(3, 7).Value = $row4[6]

# Write old row 3 data into row 4, but "Reprobadas" becomes 1
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(4, $c).Value = $row3[$c - 1]
}
$ws.Cells.Item(4, 7).Value = 1
